# "update final entrega 3"
#
# Adds the four new activity rows for the third deliverable ("Entrega 3")
# to the tracking sheet, clears the stray underline formatting that was
# left on A15, and leaves the selection on A15 (matching the author's
# final state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New activity rows (rows 16-19) ---------------------------------------
$ws.Range("A16").Value = "Rafael, Willian, Gabriel"
$ws.Range("B16").Value = "Estudo de uma estratégia para o modelo"
$ws.Range("C16").Value = "Alta"

$ws.Range("A17").Value = "Rafael"
$ws.Range("B17").Value = "Modelagem das tabelas financeiras"
$ws.Range("C17").Value = "Alta "

$ws.Range("A18").Value = "Gabriel"
$ws.Range("B18").Value = "Criação do modelo inicial"
$ws.Range("C18").Value = "Alta"

$ws.Range("A19").Value = "Willian"
$ws.Range("B19").Value = "Revisão, correção e complementação"
$ws.Range("C19").Value = "Alta"

# --- Remove the leftover underline style on A15 ---------------------------
$ws.Range("A15").Font.Underline = $false

# --- Leave the selection on A15, as in the saved file ----------------------
$ws.Range("A15").Select()
